$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").Value = "'15"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'21.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("G3").Value = "'15"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.401"
$ws.Range("D4").Style = "Normal"
$ws.Range("G4").Value = "'15"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.05847"
$ws.Range("D5").Style = "Normal"
$ws.Range("G5").Value = "'15"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'3.394"
$ws.Range("D6").Style = "Normal"
$ws.Range("G6").Value = "'15"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'6.365"
$ws.Range("D7").Style = "Normal"
$ws.Range("G7").Value = "'15"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.8145"
$ws.Range("D8").Style = "Normal"
$ws.Range("G8").Value = "'15"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'1.014"
$ws.Range("D9").Style = "Normal"
$ws.Range("G9").Value = "'15"
$ws.Range("G9").Style = "Normal"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01116"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("G10").Value = "'15"
$ws.Range("G10").Style = "Normal"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1418"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("G11").Value = "'15"
$ws.Range("G11").Style = "Normal"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03657"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value = "'15"
$ws.Range("G12").Style = "Normal"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.07414"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12MandalaExchangeTokenMDX"
$ws.Range("G13").Value = "'15"
$ws.Range("G13").Style = "Normal"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03040"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("G14").Value = "'15"
$ws.Range("G14").Style = "Normal"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'4.210"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").Value = "'15"
$ws.Range("G15").Style = "Normal"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").Value = "'0.09385"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("G16").Value = "'15"
$ws.Range("G16").Style = "Normal"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001600"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("G17").Value = "'15"
$ws.Range("G17").Style = "Normal"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04803"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("G18").Value = "'15"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'0.006009"
$ws.Range("D19").Style = "Normal"
$ws.Range("G19").Value = "'15"
$ws.Range("G19").Style = "Normal"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.004079"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("G20").Value = "'15"
$ws.Range("G20").Style = "Normal"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009973"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("G21").Value = "'15"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D22").Style = "Normal"
$ws.Range("G22").Value = "'15"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'3.694"
$ws.Range("D23").Style = "Normal"
$ws.Range("G23").Value = "'15"
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").Value = "'2.222"
$ws.Range("D24").Style = "Normal"
$ws.Range("G24").Value = "'15"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.3236"
$ws.Range("D25").Style = "Normal"
$ws.Range("G25").Value = "'15"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.1297"
$ws.Range("D26").Style = "Normal"
$ws.Range("G26").Value = "'15"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.0002502"
$ws.Range("D27").Style = "Normal"
$ws.Range("G27").Value = "'15"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'15"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'15"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'15"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'15"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'15"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'15"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'15"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'15"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'15"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'15"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'15"
$ws.Range("G38").Style = "Normal"
$ws.Range("G39").Value = "'15"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.03848"
$ws.Range("D40").Style = "Normal"
$ws.Range("G40").Value = "'15"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.006425"
$ws.Range("D41").Style = "Normal"
$ws.Range("G41").Value = "'15"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.1070"
$ws.Range("D42").Style = "Normal"
$ws.Range("G42").Value = "'15"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.002502"
$ws.Range("D43").Style = "Normal"
$ws.Range("G43").Value = "'15"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.006227"
$ws.Range("D44").Style = "Normal"
$ws.Range("G44").Value = "'15"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005627"
$ws.Range("D45").Style = "Normal"
$ws.Range("G45").Value = "'15"
$ws.Range("G45").Style = "Normal"
$ws.Range("G46").Value = "'15"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.6005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("G47").Value = "'15"
$ws.Range("G47").Style = "Normal"
$ws.Range("G48").Value = "'15"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("G49").Value = "'15"
$ws.Range("G49").Style = "Normal"
$ws.Range("G50").Value = "'15"
$ws.Range("G50").Style = "Normal"
$ws.Range("G51").Value = "'15"
$ws.Range("G51").Style = "Normal"
